$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-14 down to 3-15)
$ws.Rows("2:2").Insert()

# The insert carries blank styled placeholders into row 2 for every
# column; the new "import external data" row only has feature/owner/
# actual-time values, so drop the Estd-Time (C) and Test-Estd (E)
# placeholders entirely before (re)writing the row's real content.
$ws.Range("C2").Clear()
$ws.Range("E2").Clear()

# New feature row: "import external data"
$ws.Range("A2").Value = "import external data"
$ws.Range("B2").Value = "Alex"
$ws.Range("D2").Value = 14
$ws.Range("F2").Value = 10

# "Registration and login of users;" row picked up Actual-time figures
$ws.Range("D3").Value = 16
$ws.Range("F3").Value = 16

# "Search by product name / type;" row picked up Actual-time figures
$ws.Range("D4").Value = 16
$ws.Range("F4").Value = 16

# "Engine + core logic" row picked up Actual-time figures
$ws.Range("D10").Value = 6
$ws.Range("F10").Value = 6

# "Order model + configuration" row picked up Actual-time figures
$ws.Range("D13").Value = 8
$ws.Range("F13").Value = 8

# Restore the selected cell shown in the saved workbook
$ws.Range("F5").Select()
